# Weekly update: insert a new weekly record at row 47 ("Hortaliza, Feria
# Lagunitas de Puerto Montt - Espinaca"), pushing the existing rows 47-96
# down to 48-97. The new row repeats the constant dimensions of the table
# and carries the new week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 47, shifting rows 47:96 down to 48:97.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly observation.
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 45233
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112012
$ws.Cells.Item(47, 7).Value = "Espinaca"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 35
$ws.Cells.Item(47, 11).Value = 15000
$ws.Cells.Item(47, 12).Value = 15000
$ws.Cells.Item(47, 13).Value = 15000
$ws.Cells.Item(47, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 1500
$ws.Cells.Item(47, 17).Value = 10
$ws.Cells.Item(47, 18).Value = "Hortaliza"
